$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.022.00'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.45%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.752.12'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.96%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9971'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.33%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.39'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.17%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9971'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.32%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5105'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +8.59%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3591'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.58%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.32'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.14%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07269'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.09%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.068'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.57%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9981'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.16%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.30'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.75%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.991'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.31%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.743.96'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.50%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.842'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.68%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '87.22'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.82%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001034'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.44%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06414'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.83%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9984'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.17%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.63'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.62%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.772'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.58%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.079.76'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.37%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.27'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.77%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.039'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.89%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '152.48'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.08%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.93'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.32%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.938.25'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.26%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.250'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.94%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.35'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.42%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.047'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.67%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09618'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.81%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.590'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.13%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.406'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.53%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.05909'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.52%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02187'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.55%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '11.02'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.37%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2005'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.56%  '

$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.789'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.06%  '

$ws.Range('B40').Value = 'WEMIXTOKEN'
$ws.Range('C40').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.422'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.25%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6055'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.64%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.108'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.48%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.638'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.62%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.99'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.30%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.594'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.17%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5678'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.53%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '120.46'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.11%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.853'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.82%  '

$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06694'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.72%  '

$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.104'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.65%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9992'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.08%  '
